$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.999999841026622
$ws.Range("E2").Value = 0.999999841026622

# Row 3
$ws.Range("D3").Value = 0.9999964841144882
$ws.Range("E3").Value = 0.9999964841144882

# Row 4
$ws.Range("D4").Value = 0.5473885312129469
$ws.Range("E4").Value = 0.5473885312129469

# Row 5
$ws.Range("D5").Value = 0.9999999999993945
$ws.Range("E5").Value = 0.9999999999993945

# Row 6
$ws.Range("D6").Value = 0.986566140412361
$ws.Range("E6").Value = 0.986566140412361

# Row 7
$ws.Range("D7").Value = 0.9386490327717458
$ws.Range("E7").Value = 0.06135096722825417

# Row 8
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = 0.0000000003203158577392164
$ws.Range("E8").Value = 0.9999999996796841

# Row 9
$ws.Range("D9").Value = 0.9988106881357675
$ws.Range("E9").Value = 0.001189311864232545

# Row 10
$ws.Range("C10").Value = $false
$ws.Range("D10").Value = 0.002130186810853152
$ws.Range("E10").Value = 0.9978698131891468

# Row 11
$ws.Range("D11").Value = 0.9999999999968174
$ws.Range("E11").Value = 0.000000000003182565322390474
$ws.Range("F11").Value = 8.95258617401123
$ws.Range("G11").Value = 0.3

# Row 12
$ws.Range("D12").Value = 0.9999999963379553
$ws.Range("E12").Value = 0.9999999963379553

# Row 13
$ws.Range("D13").Value = 0.9950529196379437
$ws.Range("E13").Value = 0.9950529196379437

# Row 14
$ws.Range("D14").Value = 0.5641094717838079
$ws.Range("E14").Value = 0.5641094717838079

# Row 15
$ws.Range("D15").Value = 0.00000000000000000000001431091321125057
$ws.Range("E15").Value = 0.00000000000000000000001431091321125057

# Row 16
$ws.Range("D16").Value = 0.9994776043491816
$ws.Range("E16").Value = 0.9994776043491816

# Row 17
$ws.Range("D17").Value = 0.8761437685956361
$ws.Range("E17").Value = 0.1238562314043639

# Row 18
$ws.Range("C18").Value = $false
$ws.Range("D18").Value = 0.01977734136978154
$ws.Range("E18").Value = 0.9802226586302185

# Row 19
$ws.Range("D19").Value = 0.9998489536455246
$ws.Range("E19").Value = 0.0001510463544753815

# Row 20
$ws.Range("C20").Value = $false
$ws.Range("D20").Value = 0.0005426545148217582
$ws.Range("E20").Value = 0.9994573454851783

# Row 21
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 4.46962833404541
$ws.Range("G21").Value = 0.4
